# Apply commit: "Fruta / hortaliza, semanal"
# Insert two new weekly price rows (Primera / Segunda, fecha 2023-11-07)
# at the top of the data block (row 7/8), pushing the existing historical
# rows down by two (rows 7..91 -> 9..93). This matches the diff, where
# every existing record's row number increases by 2 and two brand new
# records appear at rows 7-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows right after the header block (before the old row 7),
# shifting all the existing data rows (old 7..91) down to 9..93.
$ws.Rows("7:8").Insert()

$ws.Cells.Item(7,1).Value = 1
$ws.Cells.Item(7,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7,3).Value = "Arica y Parinacota"
$ws.Cells.Item(7,4).Value = 45237
$ws.Cells.Item(7,5).Value = 15
$ws.Cells.Item(7,6).Value = 100112028
$ws.Cells.Item(7,7).Value = "Sandia"
$ws.Cells.Item(7,8).Value = "Sin especificar"
$ws.Cells.Item(7,9).Value = "Primera"
$ws.Cells.Item(7,10).Value = 700
$ws.Cells.Item(7,11).Value = 490
$ws.Cells.Item(7,12).Value = 500
$ws.Cells.Item(7,13).Value = 495
$ws.Cells.Item(7,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(7,15).Value = "Perú"
$ws.Cells.Item(7,16).Value = 495
$ws.Cells.Item(7,17).Value = 1
$ws.Cells.Item(7,18).Value = "Hortaliza"

$ws.Cells.Item(8,1).Value = 1
$ws.Cells.Item(8,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8,3).Value = "Arica y Parinacota"
$ws.Cells.Item(8,4).Value = 45237
$ws.Cells.Item(8,5).Value = 15
$ws.Cells.Item(8,6).Value = 100112028
$ws.Cells.Item(8,7).Value = "Sandia"
$ws.Cells.Item(8,8).Value = "Sin especificar"
$ws.Cells.Item(8,9).Value = "Segunda"
$ws.Cells.Item(8,10).Value = 550
$ws.Cells.Item(8,11).Value = 490
$ws.Cells.Item(8,12).Value = 500
$ws.Cells.Item(8,13).Value = 495
$ws.Cells.Item(8,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(8,15).Value = "Perú"
$ws.Cells.Item(8,16).Value = 495
$ws.Cells.Item(8,17).Value = 1
$ws.Cells.Item(8,18).Value = "Hortaliza"

Write-Output ("Dimension now: " + $ws.UsedRange.Address())
